$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking text values (Price / Volume columns) - force Text format to preserve exact formatting
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.89"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.150"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05760"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.276"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8497"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8582"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1384"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03417"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07075"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03236"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09359"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001528"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0005969"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005904"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.544"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.216"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3123"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.491"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04112"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004157"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03754"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003543"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009393"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.33%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.35%"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.40%"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.55%"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.05%"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "6.78%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.93%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.47%"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.18%"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "5.11%"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.17%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "13.03%"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.30%"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.12%"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.49%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.35%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.35%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.51%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.94%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.20%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.89%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.19%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.90%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.90%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-7.60%"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.80%"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.01%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-5.39%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-38.50%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-6.54%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "7.77%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-10.87%"

# Plain text values (Coin name / Link columns)
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("B16").Value = "One"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("B18").Value = "LEO"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("C16").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
